# Update countries & provincias Spain
# Refresh of COVID country stats table + re-sort of some tied rows,
# which manifests as shared-string reordering (country labels swap
# between rows) plus updated numeric values, and a refreshed timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refreshed "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 04:35"

# --- Block 1: Grecia / Honduras / Guinea / Uzbekistan / Sudan (rows 74-78) ---
# Honduras jumps ahead of Guinea & Uzbekistan with fresh case numbers;
# Guinea and Uzbekistan's rows take over the data the row above used to hold.
$ws.Range("A75").Value = "Honduras"
$ws.Range("B75").Value = 2798
$ws.Range("C75").Value = 152
$ws.Range("D75").Value = 340
$ws.Range("E75").Value = 2312
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 146

$ws.Range("A76").Value = "Guinea"
$ws.Range("B76").Value = 2796
$ws.Range("D76").Value = 1263
$ws.Range("E76").Value = 1517
$ws.Range("H76").Value = 16

$ws.Range("A77").Value = "Uzbekistan"
$ws.Range("B77").Value = 2791
$ws.Range("D77").Value = 2314
$ws.Range("E77").Value = 464
$ws.Range("H77").Value = 13

# --- Block 2: Surinam / Seychelles / Groenlandia / Montserrat (rows 208-211) ---
# Seychelles moves up right after Surinam; Groenlandia & Montserrat shift down.
$ws.Range("A209").Value = "Seychelles"

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- Block 3: Papua Nueva Guinea / San Bartolome / Bonaire / Sahara Occidental (rows 213-216) ---
# San Bartolome and Bonaire, San Eustaquio y Saba move up right after Papua
# Nueva Guinea; Sahara Occidental drops after them. All three rows carry
# identical underlying figures, so only the labels change.
$ws.Range("A214").Value = "San Bartolome"
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "Sahara Occidental"
